$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bold the "Objective N" label at the start of each "Objective N: ..."
#    paragraph, leaving the rest of the sentence in regular (non-bold) text.
#    "Objective " + single digit == 11 characters for all six objectives.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Objective *:*") {
        $start = $p.Range.Start
        $r = $d.Range($start, $start + 11)
        $r.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# 2) The paragraph "The Tracker can connect to new objects should the user
#    desire." currently has a hidden "_GoBack" bookmark splitting it into
#    two runs ("connect" / " to new objects..."). Collapse it back into a
#    single run and drop the bookmark from this location (it gets moved to
#    a new paragraph near the end of the document in step 3).
#
#    Re-assigning the exact same text is a no-op in this engine, so we first
#    push through a placeholder value and then set the real text - this
#    guarantees the run split (and the bookmark sitting between the runs)
#    is actually rebuilt as a single plain run.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "The Tracker can connect to new objects should the user desire.") {
        $rng = $p.Range
        $rng.MoveEnd(1, -1)
        $rng.Text = "__TEMP_PLACEHOLDER__"

        $rng2 = $p.Range
        $rng2.MoveEnd(1, -1)
        $rng2.Text = "The Tracker can connect to new objects should the user desire."
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark (now orphaned) as its own empty
#    paragraph. It lands right before the final, totally bare empty
#    paragraph - i.e. between the blue/accent-coloured empty paragraph and
#    the last empty paragraph of the document.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n - 1)

$bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($bookmarkXml)
